$p = $ppt.ActivePresentation

# Convert an EMU offset to the point value PowerPoint's COM layer expects.
# (914400 EMU per inch, 12700 EMU per point.) A small bias away from zero
# is added before the conversion so that this host's point->EMU
# reconversion (which truncates toward zero) lands back on the exact
# integer EMU value instead of one EMU short.
function Emu-ToPt {
    param($emu)
    $sign = 1
    if ($emu -lt 0) { $sign = -1 }
    $mag = [Math]::Abs($emu)
    return ($mag + 0.4) / 12700.0 * $sign
}

function Set-ShapePosByName {
    param($Slide, $Name, $LeftEmu, $TopEmu)
    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $sh = $Slide.Shapes.Item($i)
        if ($sh.Name -eq $Name) {
            $sh.Left = Emu-ToPt($LeftEmu)
            $sh.Top = Emu-ToPt($TopEmu)
            return
        }
    }
}

# --- Slide 1 shape repositions ---
$s1 = $p.Slides.Item(1)
Set-ShapePosByName $s1 "Rechteck 5" 0 1712384
Set-ShapePosByName $s1 "Rechteck 9" -2322 4905901

# --- Slide 3 shape repositions ---
$s3 = $p.Slides.Item(3)
Set-ShapePosByName $s3 "Rechteck 15" 0 1009709
Set-ShapePosByName $s3 "Rechteck 6"  5225415 2093595
Set-ShapePosByName $s3 "Rechteck 8"  693418 7304981
Set-ShapePosByName $s3 "Rechteck 9"  949458 7812922
Set-ShapePosByName $s3 "Rechteck 10" 5130535 7681611
Set-ShapePosByName $s3 "Rechteck 13" 949458 8360752

# --- Date placeholder text on the slide master and every slide layout ---
# (cached text of the auto-updating "datetimeFigureOut" field)
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -eq "Date Placeholder 3") {
        $sh.TextFrame.TextRange.Text = "18.09.2023"
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "05.09.2023") {
                $sh.TextFrame.TextRange.Text = "18.09.2023"
            }
        }
    }
}
